{"js": "// Each entry is [oldText, newText]; these are the exact, non-duplicated\n// strings found in the document body (the title date plus the 25\n// two-digit-by-two-digit multiplication prompts in the table).\nconst replacements = [\n  [\"2025-01-14 Tuesday\", \"2025-01-15 Wednesday\"],\n  [\"39\u00d785=\", \"34\u00d739=\"],\n  [\"92\u00d774=\", \"62\u00d766=\"],\n  [\"97\u00d747=\", \"51\u00d761=\"],\n  [\"78\u00d721=\", \"63\u00d741=\"],\n  [\"84\u00d730=\", \"76\u00d732=\"],\n  [\"95\u00d797=\", \"67\u00d726=\"],\n  [\"57\u00d719=\", \"63\u00d752=\"],\n  [\"71\u00d755=\", \"28\u00d756=\"],\n  [\"62\u00d771=\", \"37\u00d744=\"],\n  [\"99\u00d754=\", \"60\u00d723=\"],\n  [\"14\u00d750=\", \"87\u00d758=\"],\n  [\"41\u00d747=\", \"67\u00d759=\"],\n  [\"86\u00d776=\", \"86\u00d738=\"],\n  [\"82\u00d711=\", \"90\u00d753=\"],\n  [\"56\u00d763=\", \"48\u00d777=\"],\n  [\"45\u00d756=\", \"72\u00d778=\"],\n  [\"23\u00d713=\", \"43\u00d722=\"],\n  [\"99\u00d788=\", \"45\u00d799=\"],\n  [\"93\u00d754=\", \"44\u00d752=\"],\n  [\"31\u00d733=\", \"59\u00d748=\"],\n  [\"11\u00d764=\", \"33\u00d724=\"],\n  [\"35\u00d793=\", \"50\u00d797=\"],\n  [\"56\u00d734=\", \"81\u00d775=\"],\n  [\"45\u00d787=\", \"54\u00d767=\"],\n  [\"51\u00d748=\", \"79\u00d775=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Each pair is (oldText, newText); these are the exact, non-duplicated\n# strings found in the document body (the title date plus the 25\n# two-digit-by-two-digit multiplication prompts in the table).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-01-14 Tuesday', '2025-01-15 Wednesday'),\n    @('39\u00d785=', '34\u00d739='),\n    @('92\u00d774=', '62\u00d766='),\n    @('97\u00d747=', '51\u00d761='),\n    @('78\u00d721=', '63\u00d741='),\n    @('84\u00d730=', '76\u00d732='),\n    @('95\u00d797=', '67\u00d726='),\n    @('57\u00d719=', '63\u00d752='),\n    @('71\u00d755=', '28\u00d756='),\n    @('62\u00d771=', '37\u00d744='),\n    @('99\u00d754=', '60\u00d723='),\n    @('14\u00d750=', '87\u00d758='),\n    @('41\u00d747=', '67\u00d759='),\n    @('86\u00d776=', '86\u00d738='),\n    @('82\u00d711=', '90\u00d753='),\n    @('56\u00d763=', '48\u00d777='),\n    @('45\u00d756=', '72\u00d778='),\n    @('23\u00d713=', '43\u00d722='),\n    @('99\u00d788=', '45\u00d799='),\n    @('93\u00d754=', '44\u00d752='),\n    @('31\u00d733=', '59\u00d748='),\n    @('11\u00d764=', '33\u00d724='),\n    @('35\u00d793=', '50\u00d797='),\n    @('56\u00d734=', '81\u00d775='),\n    @('45\u00d787=', '54\u00d767='),\n    @('51\u00d748=', '79\u00d775='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap:=wdFindContinue, Format,\n    #          ReplaceWith, Replace:=wdReplaceAll)\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
